# Add a new slide ("环境准备" / JDK / IDE / Maven) to the 微服务架构 deck.
#
# Summary of the edit:
#   1. Duplicate slide 1 (项目实践 / 需求分析) to become the new, final slide
#      ("项目实践" / "环境准备") and give it a new text box listing the
#      environment-prep steps (JDK / IDE / Maven).
#   2. Move that duplicate to the end of the deck (slide 6).
#   3. Retitle the original slide 1's bottom caption from "需求分析" to
#      "开发流程" (it now introduces the whole 开发流程 section).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Duplicate slide 1 -> becomes the new "environment prep" slide.
# ---------------------------------------------------------------------
$origSlide = $p.Slides.Item(1)
$dupRange  = $origSlide.Duplicate()
$newSlide  = $dupRange.Item(1)

# Move the duplicate to the very end of the slide list (slide 6).
$newSlide.MoveTo($p.Slides.Count)

# ---------------------------------------------------------------------
# 2. On the new (last) slide: change the caption text and add the
#    JDK / IDE / Maven text box.
# ---------------------------------------------------------------------
$newCaption = $newSlide.Shapes.Item("Rectangle 12")
$newCaption.TextFrame.TextRange.Paragraphs(3, 1).Text = "环境准备"

$tb = $newSlide.Shapes.AddTextbox(1, 409.1051968503937, 233.64842519685038, 96.33165354330708, 72.7031496062992)
$tb.Name = "文本框 2"

$tr = $tb.TextFrame.TextRange
$tr.Text = "JDK"
$tr.InsertAfter("`rIDE")
$tr.InsertAfter("`rMaven")

$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1

$paraCount = $tb.TextFrame.TextRange.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tb.TextFrame.TextRange.Paragraphs($i, 1)
    $para.ParagraphFormat.Bullet.Type = 2
    $para.ParagraphFormat.Bullet.Font.Name = "+mj-lt"
}

# ---------------------------------------------------------------------
# 3. On the original slide 1: rename the caption to "开发流程".
# ---------------------------------------------------------------------
$origCaption = $origSlide.Shapes.Item("Rectangle 12")
$origCaption.TextFrame.TextRange.Paragraphs(3, 1).Text = "开发流程"
